$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "45.395.71"
$ws.Range("E2").Value = "  -0.11%  "
$ws.Range("D3").Value = "2.368.65"
$ws.Range("E3").Value = "  -0.43%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'312.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.89%  "
$ws.Range("D6").Value = "'107.76"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.27%  "
$ws.Range("E7").Value = "  -1.19%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "'0.609"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.50%  "
$ws.Range("D10").Value = "'40.74"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.09%  "
$ws.Range("E11").Value = "  -1.41%  "
$ws.Range("D12").Value = "'8.45"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.61%  "
$ws.Range("D13").Value = "'0.110"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.02%  "
$ws.Range("E14").Value = "  -3.92%  "
$ws.Range("D15").Value = "2.728.02"
$ws.Range("E15").Value = "  -0.37%  "
$ws.Range("D16").Value = "'15.31"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.90%  "
$ws.Range("D17").Value = "2.373.30"
$ws.Range("E17").Value = "  -0.80%  "
$ws.Range("D18").Value = "45.426.37"
$ws.Range("E18").Value = "  +0.42%  "
$ws.Range("D19").Value = "'14.11"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +8.03%  "
$ws.Range("E20").Value = "  -1.97%  "
$ws.Range("E21").Value = "  -5.53%  "
$ws.Range("D22").Value = "'73.24"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.58%  "
$ws.Range("D23").Value = "'3.53"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.82%  "
$ws.Range("D24").Value = "'259.99"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'2.38"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.58%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").Value = "'11.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.18%  "
$ws.Range("D28").Value = "'7.28"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.19%  "
$ws.Range("E29").Value = "  -1.69%  "
$ws.Range("D30").Value = "'0.0972"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.80%  "
$ws.Range("E31").Value = "  -2.86%  "
$ws.Range("D32").Value = "'36.72"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.67%  "
$ws.Range("D33").Value = "'166.26"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.11%  "
$ws.Range("D34").Value = "'2.95"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.74%  "
$ws.Range("E35").Value = "  -1.95%  "
$ws.Range("D36").Value = "'0.117"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.44%  "
$ws.Range("D37").Value = "'4.70"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.89%  "
$ws.Range("E38").Value = "  +9.18%  "
$ws.Range("E39").Value = "  +0.34%  "
$ws.Range("D40").Value = "'2.94"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.99%  "
$ws.Range("E41").Value = "  -3.63%  "
$ws.Range("D42").Value = "'98.77"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.02%  "
$ws.Range("D43").Value = "'69.91"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.61%  "
$ws.Range("D44").Value = "'0.226"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.74%  "
$ws.Range("E45").Value = "  +0.11%  "
$ws.Range("E46").Value = "  -8.48%  "
$ws.Range("D47").Value = "1.812.42"
$ws.Range("E47").Value = "  +9.55%  "
$ws.Range("D48").Value = "'83.49"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.01%  "
$ws.Range("D49").Value = "'5.81"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.82%  "
$ws.Range("D50").Value = "'111.14"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.92%  "
$ws.Range("D51").Value = "'9.21"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.56%  "
